$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 188 ("こんにちは" post) entirely; all subsequent rows shift up by one.
$ws.Rows("188:188").Delete()

# The engine round-trips previously "empty" inline-string cells (self-closed
# <is/>, i.e. no value at all) as an empty-string value on save. Restore the
# handful of cells that are genuinely blank (both pre-existing ones above the
# deleted row, and ones shifted up from below it) back to true empty cells.
$blankCells = @(
    @(47, 2),
    @(229, 3),
    @(232, 2),
    @(236, 2),
    @(282, 2),
    @(340, 2)
)
foreach ($cell in $blankCells) {
    $ws.Cells.Item($cell[0], $cell[1]).ClearContents()
}
